# "add data request handle"
#
# Inserts a new "failure_code" lookup table (FAILED_AUTH_CHECK / 身份验证错误)
# right after the existing failure_code table, pushing the pre-existing
# request_type table down by four rows (one new 3-row table + the blank
# separator row that used to sit above it).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Make room: rows 10-12 currently hold the "request_type" table, with row 13
# blank. Insert 4 rows at row 10 so that table (and its trailing blank
# separator) slides down to rows 14-17, leaving rows 10-13 free (10-12 for
# the new table, 13 stays blank just like the other table separators).
$ws.Rows("10:13").Insert()

# New table: failure_code / 值 / 含义 (mirrors the first failure_code table,
# with a new FAILED_AUTH_CHECK row appended).
$ws.Range("A10").Value = "failure_code"
$ws.Range("B10").Value = "值"
$ws.Range("C10").Value = "含义"

$ws.Range("A11").Value = "FAILED_PROTOCOL_ERR"
$ws.Range("B11").Value = "0x00000001"
$ws.Range("C11").Value = "协议错误"

$ws.Range("A12").Value = "FAILED_AUTH_CHECK"
$ws.Range("B12").Value = "0x00000002"
$ws.Range("C12").Value = "身份验证错误"

# Match the author's selection state recorded after the edit.
$ws.Range("A10:C12").Select()
